$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Rewrite the "robber caught" sentence in the Elevator Pitch blurb.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "they are caught, game is over. The aim of the ai is",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "the robber is caught, they disappear from the map and would not be able to continue stealing from the police. The ai needs",
    2
)

# ------------------------------------------------------------------
# 2) The "_GoBack" bookmark moves from the Elevator Pitch paragraph to
#    the very last (empty) paragraph of the document, right before the
#    section break. Delete the old one and recreate it there.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$lastRange = $lastPara.Range

# The target paragraph has no runs at all, and this host can't anchor a
# zero-length Range that doesn't fall inside existing run text. Work
# around that by temporarily inserting a placeholder character, anchoring
# the bookmark next to it, then deleting the placeholder again (the
# bookmark survives the deletion, collapsing back to a point).
$lastRange.InsertAfter("X")
$lastRange2 = $d.Paragraphs($lastIndex).Range
$anchor = $d.Range($lastRange2.Start, $lastRange2.Start)
$d.Bookmarks.Add("_GoBack", $anchor)
$placeholder = $d.Range($lastRange2.Start, $lastRange2.Start + 1)
$placeholder.Text = ""

# ------------------------------------------------------------------
# 3) Remove one of the redundant empty "No Spacing" paragraphs that sit
#    between the spacing-only paragraph and the "1.0 Revision History"
#    heading.
# ------------------------------------------------------------------
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "`r" -and $p.Range.ParagraphStyle.NameLocal -eq "No Spacing,SAG No Spacing") {
        $prev = $d.Paragraphs($i - 1)
        if ($prev.Range.Text -eq "`r" -and $prev.Range.ParagraphStyle.NameLocal -ne "No Spacing,SAG No Spacing") {
            $targetIndex = $i
            break
        }
    }
}
if ($targetIndex -gt 0) {
    $d.Paragraphs($targetIndex).Range.Delete()
}
